$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 468, shifting existing rows 468:485 down to 469:486.
$ws.Rows.Item(468).Insert()

# Populate the newly inserted row with the new Ajo (garlic) price record.
$ws.Cells.Item(468, 1).Value = 5
$ws.Cells.Item(468, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(468, 3).Value = "Maule"
$ws.Cells.Item(468, 4).Value = 45075
$ws.Cells.Item(468, 5).Value = 7
$ws.Cells.Item(468, 6).Value = 100112003
$ws.Cells.Item(468, 7).Value = "Ajo"
$ws.Cells.Item(468, 8).Value = "Chino"
$ws.Cells.Item(468, 9).Value = "Primera"
$ws.Cells.Item(468, 10).Value = 300
$ws.Cells.Item(468, 11).Value = 17000
$ws.Cells.Item(468, 12).Value = 17000
$ws.Cells.Item(468, 13).Value = 17000
$ws.Cells.Item(468, 14).Value = "`$/malla 10 kilos"
$ws.Cells.Item(468, 15).Value = "China"
$ws.Cells.Item(468, 16).Value = 1700
$ws.Cells.Item(468, 17).Value = 10
$ws.Cells.Item(468, 18).Value = "Hortaliza"
